# Delivery Changes: append a new row of data (row 16) below the
# existing table, matching how a user would type it in via the GUI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("2025-04-25") and D ("202") look like a date / a number, so a
# plain Value assignment would get auto-converted to a date serial / a
# numeric value. Prefix with an apostrophe (exactly like typing '202 into
# Excel) to force literal text, then clear the resulting "quote prefix"
# cell format so the cell ends up with the default (unstyled) text value,
# matching the rest of the sheet's plain text cells.
$ws.Range("A16").Value = "'2025-04-25"
$ws.Range("A16").ClearFormats()

$ws.Range("B16").Value = "Nope (S00123)"
$ws.Range("C16").Value = "This is only a test"

$ws.Range("D16").Value = "'202"
$ws.Range("D16").ClearFormats()
